$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

# --- Append new trade row (row 15): Trade #16 closed info ---
$row = 15

$ws.Cells.Item($row, 1).Value = 16
# Column B holds a literal "YYYY-MM-DD" text value (not a real date). Writing it
# directly would make Excel auto-convert it to a date serial, so instead we
# build the literal string via a text formula on a scratch cell and paste the
# computed value back in - this preserves it as plain text.
$ws.Range("ZZ1").Formula = '="2026-02-16"'
$ws.Range("ZZ1").Copy()
$ws.Range("B" + $row).PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Cells.Item($row, 3).Value = "21:24:38"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"
$ws.Cells.Item($row, 6).Value = 69303.36500000001
# Column G (Exit Price) stays blank - trade is still OPEN.
$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.75
$ws.Cells.Item($row, 12).Value = "Coinbase leading with -0.081% move"
# Column M (Exit Reason) stays blank - trade is still OPEN.
$ws.Cells.Item($row, 14).Value = 0

# --- Widen the "Entry Reason" column (L) from 35 to 36 characters ---
$ws.Columns.Item(12).ColumnWidth = 36 - 5/6
